# Applies the "Automatic update of files" edit: the artfynd export swapped
# the reported order of two species-pairs of records (row 2 <-> row 4, and
# row 3 <-> row 5), including moving the "Publik kommentar" free-text note
# that travels with row 2's sighting over to row 4.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (becomes the former row 4 sighting; loses the AC "Publik kommentar") ---
$ws.Range("A2").Value = 112491413
$ws.Range("B2").Value = 90830
$ws.Range("E2").Value = 2059
$ws.Range("F2").Value = "Skrovlig taggsvamp"
$ws.Range("G2").Value = "Hydnellum scabrosum"
$ws.Range("H2").Value = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Range("P2").Value = "Älggropsröset, Vrm"
$ws.Range("Q2").Value = 356670
$ws.Range("R2").Value = 6742658
$ws.Range("AC2").ClearContents()

# --- Row 3 (becomes the former row 5 sighting) ---
$ws.Range("A3").Value = 112491434
$ws.Range("B3").Value = 90857
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 5448
$ws.Range("F3").Value = "Svartvit taggsvamp"
$ws.Range("G3").Value = "Phellodon connatus"
$ws.Range("H3").Value = "(Schultz) nom.prov"

# --- Row 4 (becomes the former row 2 sighting; gains the AC "Publik kommentar") ---
$ws.Range("A4").Value = 112491312
$ws.Range("B4").Value = 73834
$ws.Range("E4").Value = 6440
$ws.Range("F4").Value = "Vitgrynig nållav"
$ws.Range("G4").Value = "Chaenotheca subroscida"
$ws.Range("H4").Value = "(Eitner) Zahlbr."
$ws.Range("P4").Value = "Kobäcken, Vrm"
$ws.Range("Q4").Value = 356888
$ws.Range("R4").Value = 6742975
$ws.Range("AC4").Value = "På gammal senvuxen gran"

# --- Row 5 (becomes the former row 3 sighting) ---
$ws.Range("A5").Value = 112491430
$ws.Range("B5").Value = 90814
$ws.Range("D5").Value = "LC"
$ws.Range("E5").Value = 4364
$ws.Range("F5").Value = "Dropptaggsvamp"
$ws.Range("G5").Value = "Hydnellum ferrugineum"
$ws.Range("H5").Value = "(Fr.:Fr.) P. Karst."
